$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    "C2" = 4.913690000476793
    "D2" = 4.250606414976814
    "E2" = 10.09188513958596
    "F2" = 53.72696604482505
    "G2" = 3.766331770397462
    "J2" = 9.793859989880813
    "K2" = 24.23903104987336
    "M2" = 22.31489567950007
    "N2" = 22.0636284640353
    "C3" = 4.761077490840225
    "D3" = 4.233886787091286
    "E3" = 10.10962414166382
    "F3" = 53.56573533160486
    "G3" = 3.771080155468673
    "J3" = 9.816219736785712
    "K3" = 23.94251618050755
    "M3" = 22.22319591387077
    "N3" = 22.12693675646465
    "C4" = 4.666623736927143
    "D4" = 4.223408869024448
    "E4" = 10.12148923093118
    "F4" = 53.48050521123295
    "G4" = 3.774143060359801
    "J4" = 9.830944162554639
    "K4" = 23.76608522142176
    "M4" = 22.17228974937879
    "N4" = 22.16782978745203
    "C5" = 4.628012139395276
    "D5" = 4.219085179588879
    "E5" = 10.12656953481512
    "F5" = 53.44924345338561
    "G5" = 3.775428437910564
    "J5" = 9.837195137229294
    "K5" = 23.69568559819332
    "M5" = 22.15291530638197
    "N5" = 22.1850025741693
    "C6" = 4.621595456038571
    "D6" = 4.218363965723253
    "E6" = 10.12742793702423
    "F6" = 53.44426222586164
    "G6" = 3.775644126486571
    "J6" = 9.838248255553859
    "K6" = 23.68408855273181
    "M6" = 22.14978129508069
    "N6" = 22.18788482427315
    "C7" = 4.666103405062982
    "D7" = 4.223350777009221
    "E7" = 10.12155675234485
    "F7" = 53.48006954504118
    "G7" = 3.774160244499749
    "J7" = 9.831027450045344
    "K7" = 23.76512961911343
    "M7" = 22.17202289510202
    "N7" = 22.16805932623954
    "C8" = 4.861263883446467
    "D8" = 4.244885158430522
    "E8" = 10.09779981350146
    "F8" = 53.66852286096427
    "G8" = 3.76793852184203
    "J8" = 9.801363226879227
    "K8" = 24.13567410896959
    "M8" = 22.28216554352808
    "N8" = 22.08503773845074
    "C9" = 5.235279909434019
    "D9" = 4.285457399882805
    "E9" = 10.05891541343848
    "F9" = 54.14675676735467
    "G9" = 3.756899821110417
    "J9" = 9.751075287851375
    "K9" = 24.90325205544843
    "M9" = 22.54033845147561
    "N9" = 21.93826109186915
    "C10" = 5.501411001707773
    "D10" = 4.314292640577568
    "E10" = 10.03501669476631
    "F10" = 54.56338292459083
    "G10" = 3.74948790024233
    "J10" = 9.718914488851894
    "K10" = 25.48700268046352
    "M10" = 22.75475891359966
    "N10" = 21.8401790754977
    "C11" = 5.619986257944064
    "D11" = 4.327207154887929
    "E11" = 10.0251529862433
    "F11" = 54.76679851335318
    "G11" = 3.746265442116415
    "J11" = 9.705318722017015
    "K11" = 25.75575204376558
    "M11" = 22.85743575814604
    "N11" = 21.79767429886429
    "C12" = 5.664484962980232
    "D12" = 4.332068981047115
    "E12" = 10.02156235238552
    "F12" = 54.84579333346834
    "G12" = 3.745066474531093
    "J12" = 9.700318791489881
    "K12" = 25.85788703389455
    "M12" = 22.89703274266678
    "N12" = 21.78188268069226
    "C13" = 5.654919972875859
    "D13" = 4.331023164105757
    "E13" = 10.02232923763905
    "F13" = 54.82869349977783
    "G13" = 3.745323748489116
    "J13" = 9.701389015631815
    "K13" = 25.83587567573252
    "M13" = 22.88847338245444
    "N13" = 21.78527016419927
    "C14" = 5.623655515957146
    "D14" = 4.327607710749027
    "E14" = 10.02485468816654
    "F14" = 54.77325826209407
    "G14" = 3.746166376196153
    "J14" = 9.704904400119158
    "K14" = 25.76414796419056
    "M14" = 22.86067922295789
    "N14" = 21.79636901901079
    "C15" = 5.604451365308312
    "D15" = 4.325511935781762
    "E15" = 10.02642041080072
    "F15" = 54.73955757705596
    "G15" = 3.746685280206227
    "J15" = 9.70707700540088
    "K15" = 25.72025747501105
    "M15" = 22.84374696982597
    "N15" = 21.80320698506027
    "C16" = 5.493607768657553
    "D16" = 4.313444612351144
    "E16" = 10.03568155925021
    "F16" = 54.55036612126969
    "G16" = 3.749701485627103
    "J16" = 9.719823796613973
    "K16" = 25.46949553748966
    "M16" = 22.74815014791763
    "N16" = 21.84299939698698
    "C17" = 5.424937613956811
    "D17" = 4.305990147207281
    "E17" = 10.04162084494986
    "F17" = 54.43783954567294
    "G17" = 3.751589948257134
    "J17" = 9.727908282714171
    "K17" = 25.31641286955775
    "M17" = 22.69080404772764
    "N17" = 21.86795208547366
    "C18" = 5.385208472174677
    "D18" = 4.301683417055528
    "E18" = 10.04513186094279
    "F18" = 54.37442678189001
    "G18" = 3.752690200590159
    "J18" = 9.732655636721214
    "K18" = 25.22867041999407
    "M18" = 22.6583049701739
    "N18" = 21.88250310769888
    "C19" = 5.371718496127227
    "D19" = 4.300221932099054
    "E19" = 10.04633694135071
    "F19" = 54.3531821306707
    "G19" = 3.753065146459279
    "J19" = 9.734279740448681
    "K19" = 25.19901789552094
    "M19" = 22.6473853003602
    "N19" = 21.88746399074402
    "C20" = 5.432272010898012
    "D20" = 4.306785662576161
    "E20" = 10.04097877939318
    "F20" = 54.44968287662309
    "G20" = 3.751387464310291
    "J20" = 9.727037599601426
    "K20" = 25.33267773105113
    "M20" = 22.69685861727137
    "N20" = 21.86527524336997
    "C21" = 5.632849924290055
    "D21" = 4.328611684637275
    "E21" = 10.02410898259585
    "F21" = 54.78948785505917
    "G21" = 3.745918298981973
    "J21" = 9.703867818391771
    "K21" = 25.78520696135773
    "M21" = 22.86882381314174
    "N21" = 21.79310076400989
    "C22" = 5.761571445792005
    "D22" = 4.342709682253185
    "E22" = 10.01392589012266
    "F22" = 55.02301081206592
    "G22" = 3.742468006430158
    "J22" = 9.689590439055268
    "K22" = 26.08304811087873
    "M22" = 22.98537187221203
    "N22" = 21.74770271634655
    "C23" = 5.693100751015079
    "D23" = 4.335200378469827
    "E23" = 10.01928386076176
    "F23" = 54.89733974789046
    "G23" = 3.744298186849687
    "J23" = 9.697131439143382
    "K23" = 25.92392461004582
    "M23" = 22.92279538338115
    "N23" = 21.77177029975314
    "C24" = 5.428956904489306
    "D24" = 4.306426075712056
    "E24" = 10.0412687568071
    "F24" = 54.44432451660776
    "G24" = 3.75147896197071
    "J24" = 9.727430925574874
    "K24" = 25.32532354892247
    "M24" = 22.69411988083036
    "N24" = 21.86648480395364
    "C25" = 5.135394539148786
    "D25" = 4.274655080566814
    "E25" = 10.06861274137934
    "F25" = 54.00583310636294
    "G25" = 3.759762730847575
    "J25" = 9.763837643098928
    "K25" = 24.69174482704283
    "M25" = 22.46607399033259
    "N25" = 21.9762539692537
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
